$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (shifts existing rows 24-40 down to 25-41)
$ws.Rows("24:24").Insert()

# Populate the new row 24 with this week's record
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 44651
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 100112030
$ws.Range("G24").Value = "Poroto granado"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 28000
$ws.Range("L24").Value = 30000
$ws.Range("M24").Value = 28960
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 1158
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
